$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hashcode.csv")

$ws.Cells.Item(9, 2).Value = "b38f934c02d047a2ada11101a82c1f39"
$ws.Cells.Item(17, 2).Value = "07256692167359f375548b4159378639"
$ws.Cells.Item(44, 2).Value = "cad1b42e28fd98d0f49671f96c3de57e"
$ws.Cells.Item(89, 2).Value = "8389652cded787e07292c43f522ca13f"
$ws.Cells.Item(94, 2).Value = "44213aeeab26b84a909d27da8747f1dd"
$ws.Cells.Item(95, 2).Value = "62fb3a25e5eb73fa548e78df049eeae4"
$ws.Cells.Item(98, 2).Value = "7e28e709da59e3fc566edfc13a487028"
$ws.Cells.Item(99, 2).Value = "a903b995188f2ddede8f2f45a9506657"
$ws.Cells.Item(110, 2).Value = "7363d9afdc32195301b0eff7cd8ac049"
$ws.Cells.Item(115, 2).Value = "78fb34603fc974bb8815be6ff28d67f3"
$ws.Cells.Item(136, 2).Value = "96cc8ec8396de1965e96bbd5ad31232e"
$ws.Cells.Item(154, 2).Value = "710498cd1c9d97b7b820546131b3b3f1"
$ws.Cells.Item(159, 2).Value = "e156324346350ccc26b299ced439b2e5"
$ws.Cells.Item(160, 2).Value = "7cbf8ffcf68818fdb9ee36505a027f56"
$ws.Cells.Item(168, 2).Value = "d5d7c02dde683b92ee85060ddc3fd8c6"
$ws.Cells.Item(169, 2).Value = "574e5e040d38a98dd0601262e0d8c5ee"
$ws.Cells.Item(183, 2).Value = "b3a5b41de62bc70708855999dc05272a"
$ws.Cells.Item(200, 2).Value = "167b2fa8a52251f81750b9c2cb5d4eea"
$ws.Cells.Item(222, 2).Value = "62583ae869ae3960864909d2c138818f"
$ws.Cells.Item(227, 2).Value = "2d01a5278488f10b9f5dd5e43c9859b6"
$ws.Cells.Item(228, 2).Value = "64b0b49079d4fafbf463562b0ce5c243"
$ws.Cells.Item(229, 2).Value = "8f631ee40d39ff576db24dcf77081725"
$ws.Cells.Item(232, 2).Value = "c7017acfe56676dd01830aabf3c16619"
$ws.Cells.Item(246, 2).Value = "5a276b413bae9d54af1af76aa1369b6d"
$ws.Cells.Item(276, 2).Value = "5cf49fcb7b42f8c257473277b881c555"
$ws.Cells.Item(278, 2).Value = "7945392d2c0a38ccd06a83268ca681d3"
$ws.Cells.Item(281, 2).Value = "785770d2c4d28e162a36d58ef4a28c59"
$ws.Cells.Item(335, 2).Value = "50b58f4097323141d0561c5f09b1d665"
$ws.Cells.Item(339, 2).Value = "06373edd62c8e2b9ae8898478ade4978"
$ws.Cells.Item(411, 2).Value = "c2f5cdcfa1bdf05b5752b8b0e460e991"
$ws.Cells.Item(420, 2).Value = "bf3569543f5afe0bd329968445d710df"
$ws.Cells.Item(448, 2).Value = "a936ee92276cb2a0337d96d64f3c12bc"
$ws.Cells.Item(464, 2).Value = "bd4dd24afc8c843a92b91727d260b2e9"
$ws.Cells.Item(483, 2).Value = "7db025c699f5ae5fc290487270fbbc2d"
$ws.Cells.Item(506, 2).Value = "ea1655feed1aab34539739f88df3e2a3"
$ws.Cells.Item(507, 2).Value = "85f4d1013fdee287aa9ccfbbb2ed9b2e"
$ws.Cells.Item(508, 2).Value = "a601a1dda8c8c65aaa8e64d71cc76a02"
$ws.Cells.Item(523, 2).Value = "ce7b5d6b92bbdaa38a3432c113352a67"
$ws.Cells.Item(524, 2).Value = "3fcdff4b97f2ae0d13fa819c55d00493"
$ws.Cells.Item(555, 2).Value = "1bed013139733d869b685643fab09c7d"
$ws.Cells.Item(561, 2).Value = "b97d197b9e2546ccd903d181ecd3d671"
$ws.Cells.Item(574, 2).Value = "5e4a97d95670e8025a08c9355ea3e2d5"
$ws.Cells.Item(580, 2).Value = "fa0233183a94dd823d1a0c00a9af25d2"
$ws.Cells.Item(592, 2).Value = "beba183c47427ca50cbc89fae768a4b1"
$ws.Cells.Item(600, 2).Value = "98a7a4c7e45a4c7f13b04e8c8f695464"
$ws.Cells.Item(624, 2).Value = "8087a7ff768fba1c6fb773965ee470d9"
$ws.Cells.Item(626, 2).Value = "cdeec3a4e361cc7e3e633c7a2be1280d"
$ws.Cells.Item(635, 2).Value = "17f107c3ec809afa64d7dd72684ac46a"
$ws.Cells.Item(708, 2).Value = "4570d27e6a8bc72da567b9b4478ea4f2"
$ws.Cells.Item(723, 2).Value = "356ca7a6a0143f6e4c614d0549b08df8"
$ws.Cells.Item(764, 2).Value = "14f35261c4878b68b44382223738f324"
$ws.Cells.Item(769, 2).Value = "663018af2185307a43d8ae8065e375d8"
$ws.Cells.Item(776, 2).Value = "6b9cbcf17b85706c642899379d6e2c7a"
$ws.Cells.Item(794, 2).Value = "9694dfb366b746aee7c296369334548f"
$ws.Cells.Item(824, 2).Value = "5fd89c9a1abf163cdbf6bb742d7a6c83"
$ws.Cells.Item(827, 2).Value = "5c5abd9c802a6043120d7bf33bf913e3"
$ws.Cells.Item(833, 2).Value = "0f0cb957d8aa28dd130f051d61a7ee03"
$ws.Cells.Item(835, 2).Value = "364db6aaa975fb77025f3456cfdd9b5f"
$ws.Cells.Item(838, 2).Value = "d36af27cf9b1b7bb03d6bfa34b32b383"
$ws.Cells.Item(843, 2).Value = "8c89d3b3db4666da35294c550d2fc0f0"
$ws.Cells.Item(863, 2).Value = "78076b5df5127d06a2cb137dd62cf420"
$ws.Cells.Item(877, 2).Value = "8ee2fe4e17740ed7ea71569f13f1026d"
$ws.Cells.Item(913, 2).Value = "de22376391ec9c9be46710132f4e0d0e"
$ws.Cells.Item(937, 2).Value = "47e21db15368a8cd4bcbae1681190a38"
